# Update the run_command_4 row (row 2) of the parameter table on Sheet1
# with the new values for the sound-effect run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "EURJPY=X"   # ASSET
$ws.Range("B2").Value = "Custom"     # TYPE
$ws.Range("D2").Value = 5            # FUTURE
$ws.Range("F2").Value = "AV"         # SOURCE
$ws.Range("G2").Value = "1h"         # INTERVAL

# Columns A:B grew to fit the new, longer asset/type text.
$ws.Range("A:B").EntireColumn.AutoFit()

# Leave the selection on the cell that was last edited/reviewed.
$ws.Range("G2").Select()
